$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 649.5
$ws.Range("J9").Value = 666.6667
$ws.Range("L9").Value = 666.6667
$ws.Range("N9").Value = -1004.6667
$ws.Range("H19").Value = 5079.316
$ws.Range("J19").Value = 5787.0713
$ws.Range("L19").Value = 5787.0713
$ws.Range("N19").Value = -6137.0713
$ws.Range("H33").Value = 380.72223
$ws.Range("J33").Value = 1197
$ws.Range("L33").Value = 1197
$ws.Range("N33").Value = -1655
$ws.Range("H58").Value = 3244
$ws.Range("I58").Value = 584.125
$ws.Range("K58").Value = 1752.375
$ws.Range("M58").Value = -1602.375
$ws.Range("H80").Value = 4585.8213
$ws.Range("I80").Value = 454.45456
$ws.Range("J80").Value = 7259.0586
$ws.Range("K80").Value = 1363.36368
$ws.Range("L80").Value = 21777.1758
$ws.Range("M80").Value = -365.3636799999999
$ws.Range("N80").Value = -23773.1758
$ws.Range("H82").Value = 1000
$ws.Range("I82").Value = 1000
$ws.Range("K82").Value = 3000
$ws.Range("M82").Value = -2594
$ws.Range("H83").Value = 4585.8213
$ws.Range("I83").Value = 454.45456
$ws.Range("J83").Value = 7259.0586
$ws.Range("K83").Value = 4090.09104
$ws.Range("L83").Value = 65331.52740000001
$ws.Range("M83").Value = 901.9089599999998
$ws.Range("N83").Value = -75315.52740000001
$ws.Range("H85").Value = 1000
$ws.Range("I85").Value = 1000
$ws.Range("K85").Value = 3000
$ws.Range("M85").Value = -1596
$ws.Range("H99").Value = 369.54544
$ws.Range("I99").Value = 295.75
$ws.Range("K99").Value = 887.25
$ws.Range("M99").Value = 610.75
$ws.Range("H101").Value = 459.5
$ws.Range("I101").Value = 472.33334
$ws.Range("J101").Value = 446.66666
$ws.Range("K101").Value = 1417.00002
$ws.Range("L101").Value = 1339.99998
$ws.Range("M101").Value = 204.9999800000001
$ws.Range("N101").Value = -4583.999980000001
$ws.Range("H113").Value = 3678.5
$ws.Range("I113").Value = 3330.75
$ws.Range("J113").Value = 4026.25
$ws.Range("K113").Value = 3330.75
$ws.Range("L113").Value = 4026.25
$ws.Range("M113").Value = -76.75
$ws.Range("N113").Value = -10534.25
$ws.Range("H118").Value = 765.63635
$ws.Range("I118").Value = 712.2
$ws.Range("K118").Value = 2136.6
$ws.Range("M118").Value = -479.6000000000004
$ws.Range("H135").Value = 22291.334
$ws.Range("I135").Value = 897.9
$ws.Range("K135").Value = 8081.099999999999
$ws.Range("M135").Value = -5546.099999999999
$ws.Range("H137").Value = 13215.793
$ws.Range("I137").Value = 28223
$ws.Range("K137").Value = 84669
$ws.Range("M137").Value = -82119
$ws.Range("H138").Value = 18745.525
$ws.Range("I138").Value = 1777.1904
$ws.Range("J138").Value = 56254.473
$ws.Range("K138").Value = 5331.5712
$ws.Range("L138").Value = 168763.419
$ws.Range("M138").Value = -191.5712000000003
$ws.Range("N138").Value = -179043.419

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 15000
$ws.Range("I37").Value = 15000
$ws.Range("K37").Value = 15000
$ws.Range("M37").Value = -14727
$ws.Range("H45").Value = 3966.6924
$ws.Range("I45").Value = 2641.889
$ws.Range("K45").Value = 2641.889
$ws.Range("M45").Value = -2264.889
$ws.Range("H97").Value = 1312.258
$ws.Range("I97").Value = 994.6957
$ws.Range("J97").Value = 2225.25
$ws.Range("K97").Value = 994.6957
$ws.Range("L97").Value = 2225.25
$ws.Range("M97").Value = -498.6957
$ws.Range("N97").Value = -3217.25
$ws.Range("H102").Value = 3801.125
$ws.Range("I102").Value = 3987.2666
$ws.Range("K102").Value = 3987.2666
$ws.Range("M102").Value = -2365.2666
$ws.Range("H135").Value = 87714.5
$ws.Range("J135").Value = 87714.5
$ws.Range("L135").Value = 87714.5
$ws.Range("N135").Value = -97854.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 4073
$ws.Range("I8").Value = 5999.5
$ws.Range("J8").Value = 220
$ws.Range("K8").Value = 5999.5
$ws.Range("L8").Value = 220
$ws.Range("M8").Value = -5859.5
$ws.Range("N8").Value = -500
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("H22").Value = 991.1111
$ws.Range("I22").Value = 503.8
$ws.Range("K22").Value = 503.8
$ws.Range("M22").Value = -153.8
$ws.Range("H86").Value = 48467.895
$ws.Range("I86").Value = 69353.82000000001
$ws.Range("K86").Value = 69353.82000000001
$ws.Range("M86").Value = -68230.82000000001
$ws.Range("H89").Value = 48467.895
$ws.Range("I89").Value = 69353.82000000001
$ws.Range("K89").Value = 346769.1
$ws.Range("M89").Value = -341153.1
$ws.Range("H94").Value = 2309.75
$ws.Range("I94").Value = 2399.5
$ws.Range("J94").Value = 2279.8333
$ws.Range("K94").Value = 2399.5
$ws.Range("L94").Value = 2279.8333
$ws.Range("M94").Value = -1948.5
$ws.Range("N94").Value = -3181.8333
$ws.Range("H134").Value = 2185.8667
$ws.Range("I134").Value = 1679
$ws.Range("K134").Value = 5037
$ws.Range("M134").Value = -2502
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 9250.333000000001
$ws.Range("I34").Value = 1434
$ws.Range("J34").Value = 17066.666
$ws.Range("K34").Value = 4302
$ws.Range("L34").Value = 51199.99800000001
$ws.Range("M34").Value = -4218
$ws.Range("N34").Value = -51367.99800000001
$ws.Range("H39").Value = 10501
$ws.Range("J39").Value = 10501
$ws.Range("L39").Value = 31503
$ws.Range("N39").Value = -32091
$ws.Range("H55").Value = 907.1429000000001
$ws.Range("J55").Value = 1000
$ws.Range("L55").Value = 3000
$ws.Range("N55").Value = -3354
$ws.Range("H61").Value = 92
$ws.Range("J61").Value = 100
$ws.Range("L61").Value = 300
$ws.Range("N61").Value = -730
$ws.Range("H97").Value = 1221.12
$ws.Range("J97").Value = 1359.1052
$ws.Range("L97").Value = 4077.3156
$ws.Range("N97").Value = -5069.3156
$ws.Range("H141").Value = 6491.4
$ws.Range("I141").Value = 5630.5713
$ws.Range("J141").Value = 8500
$ws.Range("K141").Value = 16891.7139
$ws.Range("L141").Value = 25500
$ws.Range("M141").Value = -11711.7139
$ws.Range("N141").Value = -35860

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 15098.333
$ws.Range("I102").Value = 21127.807
$ws.Range("J102").Value = 1747.3572
$ws.Range("K102").Value = 21127.807
$ws.Range("L102").Value = 1747.3572
$ws.Range("M102").Value = -19505.807
$ws.Range("N102").Value = -4991.3572
$ws.Range("H126").Value = 3074.9473
$ws.Range("I126").Value = 2042.8
$ws.Range("K126").Value = 6128.4
$ws.Range("M126").Value = -3658.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3799.4
$ws.Range("I7").Value = 2998
$ws.Range("K7").Value = 2998
$ws.Range("M7").Value = -2886
$ws.Range("H126").Value = 3799.4
$ws.Range("I126").Value = 2998
$ws.Range("K126").Value = 8994
$ws.Range("M126").Value = -6524
$ws.Range("H132").Value = 2747.9333
$ws.Range("I132").Value = 2479.2856
$ws.Range("J132").Value = 3688.2
$ws.Range("K132").Value = 7437.8568
$ws.Range("L132").Value = 11064.6
$ws.Range("M132").Value = -4907.8568
$ws.Range("N132").Value = -16124.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 29108.182
$ws.Range("I136").Value = 32835.95
$ws.Range("J136").Value = 5499
$ws.Range("K136").Value = 98507.84999999999
$ws.Range("L136").Value = 16497
$ws.Range("M136").Value = -95957.84999999999
$ws.Range("N136").Value = -21597
$ws.Range("H141").Value = 87994.14
$ws.Range("J141").Value = 90993.164
$ws.Range("L141").Value = 90993.164
$ws.Range("N141").Value = -101353.164
